$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Rewrite rows 6 through 43 of the GHGs sheet with the updated
# Tech/Emission/EmissionActivityRatio data (adds CO2PAN rows alongside
# existing CO2CRI rows, sorted by Tech).
$ws1.Cells.Item(6, 1).Value = 1
$ws1.Cells.Item(6, 2).Value = "PWRCCGPANXX00"
$ws1.Cells.Item(6, 3).Value = "CO2PAN"
$ws1.Cells.Item(6, 4).Value = 0.1232
$ws1.Cells.Item(6, 5).Value = "MT"

$ws1.Cells.Item(7, 1).Value = 2
$ws1.Cells.Item(7, 2).Value = "PWRCCGPANXX00"
$ws1.Cells.Item(7, 3).Value = "CO2PAN"
$ws1.Cells.Item(7, 4).Value = 0.1232
$ws1.Cells.Item(7, 5).Value = "MT"

$ws1.Cells.Item(8, 1).Value = 1
$ws1.Cells.Item(8, 2).Value = "PWRCCGPANXX01"
$ws1.Cells.Item(8, 3).Value = "CO2PAN"
$ws1.Cells.Item(8, 4).Value = 0.1006
$ws1.Cells.Item(8, 5).Value = "MT"

$ws1.Cells.Item(9, 1).Value = 2
$ws1.Cells.Item(9, 2).Value = "PWRCCGPANXX01"
$ws1.Cells.Item(9, 3).Value = "CO2PAN"
$ws1.Cells.Item(9, 4).Value = 0.1006
$ws1.Cells.Item(9, 5).Value = "MT"

$ws1.Cells.Item(10, 1).Value = 1
$ws1.Cells.Item(10, 2).Value = "PWRCCSCRIXX01"
$ws1.Cells.Item(10, 3).Value = "CO2CRI"
$ws1.Cells.Item(10, 4).Value = 0.0273
$ws1.Cells.Item(10, 5).Value = "MT"

$ws1.Cells.Item(11, 1).Value = 2
$ws1.Cells.Item(11, 2).Value = "PWRCCSCRIXX01"
$ws1.Cells.Item(11, 3).Value = "CO2CRI"
$ws1.Cells.Item(11, 4).Value = 0.0273
$ws1.Cells.Item(11, 5).Value = "MT"

$ws1.Cells.Item(12, 1).Value = 1
$ws1.Cells.Item(12, 2).Value = "PWRCCSPANXX01"
$ws1.Cells.Item(12, 3).Value = "CO2PAN"
$ws1.Cells.Item(12, 4).Value = 0.0273
$ws1.Cells.Item(12, 5).Value = "MT"

$ws1.Cells.Item(13, 1).Value = 2
$ws1.Cells.Item(13, 2).Value = "PWRCCSPANXX01"
$ws1.Cells.Item(13, 3).Value = "CO2PAN"
$ws1.Cells.Item(13, 4).Value = 0.0273
$ws1.Cells.Item(13, 5).Value = "MT"

$ws1.Cells.Item(14, 1).Value = 1
$ws1.Cells.Item(14, 2).Value = "PWRCOACRIXX01"
$ws1.Cells.Item(14, 3).Value = "CO2CRI"
$ws1.Cells.Item(14, 4).Value = 0.2724
$ws1.Cells.Item(14, 5).Value = "MT"

$ws1.Cells.Item(15, 1).Value = 2
$ws1.Cells.Item(15, 2).Value = "PWRCOACRIXX01"
$ws1.Cells.Item(15, 3).Value = "CO2CRI"
$ws1.Cells.Item(15, 4).Value = 0.2724
$ws1.Cells.Item(15, 5).Value = "MT"

$ws1.Cells.Item(16, 1).Value = 1
$ws1.Cells.Item(16, 2).Value = "PWRCOAPANXX01"
$ws1.Cells.Item(16, 3).Value = "CO2PAN"
$ws1.Cells.Item(16, 4).Value = 0.2771
$ws1.Cells.Item(16, 5).Value = "MT"

$ws1.Cells.Item(17, 1).Value = 2
$ws1.Cells.Item(17, 2).Value = "PWRCOAPANXX01"
$ws1.Cells.Item(17, 3).Value = "CO2PAN"
$ws1.Cells.Item(17, 4).Value = 0.2771
$ws1.Cells.Item(17, 5).Value = "MT"

$ws1.Cells.Item(18, 1).Value = 1
$ws1.Cells.Item(18, 2).Value = "PWRCOGCRIXX01"
$ws1.Cells.Item(18, 3).Value = "CO2CRI"
$ws1.Cells.Item(18, 4).Value = 0.2724
$ws1.Cells.Item(18, 5).Value = "MT"

$ws1.Cells.Item(19, 1).Value = 2
$ws1.Cells.Item(19, 2).Value = "PWRCOGCRIXX01"
$ws1.Cells.Item(19, 3).Value = "CO2CRI"
$ws1.Cells.Item(19, 4).Value = 0.2724
$ws1.Cells.Item(19, 5).Value = "MT"

$ws1.Cells.Item(20, 1).Value = 1
$ws1.Cells.Item(20, 2).Value = "PWRCOGPANXX01"
$ws1.Cells.Item(20, 3).Value = "CO2PAN"
$ws1.Cells.Item(20, 4).Value = 0.2724
$ws1.Cells.Item(20, 5).Value = "MT"

$ws1.Cells.Item(21, 1).Value = 2
$ws1.Cells.Item(21, 2).Value = "PWRCOGPANXX01"
$ws1.Cells.Item(21, 3).Value = "CO2PAN"
$ws1.Cells.Item(21, 4).Value = 0.2724
$ws1.Cells.Item(21, 5).Value = "MT"

$ws1.Cells.Item(22, 1).Value = 1
$ws1.Cells.Item(22, 2).Value = "PWROCGCRIXX00"
$ws1.Cells.Item(22, 3).Value = "CO2CRI"
$ws1.Cells.Item(22, 4).Value = 0.1504
$ws1.Cells.Item(22, 5).Value = "MT"

$ws1.Cells.Item(23, 1).Value = 2
$ws1.Cells.Item(23, 2).Value = "PWROCGCRIXX00"
$ws1.Cells.Item(23, 3).Value = "CO2CRI"
$ws1.Cells.Item(23, 4).Value = 0.1504
$ws1.Cells.Item(23, 5).Value = "MT"

$ws1.Cells.Item(24, 1).Value = 1
$ws1.Cells.Item(24, 2).Value = "PWROCGCRIXX01"
$ws1.Cells.Item(24, 3).Value = "CO2CRI"
$ws1.Cells.Item(24, 4).Value = 0.1437
$ws1.Cells.Item(24, 5).Value = "MT"

$ws1.Cells.Item(25, 1).Value = 2
$ws1.Cells.Item(25, 2).Value = "PWROCGCRIXX01"
$ws1.Cells.Item(25, 3).Value = "CO2CRI"
$ws1.Cells.Item(25, 4).Value = 0.1437
$ws1.Cells.Item(25, 5).Value = "MT"

$ws1.Cells.Item(26, 1).Value = 1
$ws1.Cells.Item(26, 2).Value = "PWROCGPANXX00"
$ws1.Cells.Item(26, 3).Value = "CO2PAN"
$ws1.Cells.Item(26, 4).Value = 0.1504
$ws1.Cells.Item(26, 5).Value = "MT"

$ws1.Cells.Item(27, 1).Value = 2
$ws1.Cells.Item(27, 2).Value = "PWROCGPANXX00"
$ws1.Cells.Item(27, 3).Value = "CO2PAN"
$ws1.Cells.Item(27, 4).Value = 0.1504
$ws1.Cells.Item(27, 5).Value = "MT"

$ws1.Cells.Item(28, 1).Value = 1
$ws1.Cells.Item(28, 2).Value = "PWROCGPANXX01"
$ws1.Cells.Item(28, 3).Value = "CO2PAN"
$ws1.Cells.Item(28, 4).Value = 0.1437
$ws1.Cells.Item(28, 5).Value = "MT"

$ws1.Cells.Item(29, 1).Value = 2
$ws1.Cells.Item(29, 2).Value = "PWROCGPANXX01"
$ws1.Cells.Item(29, 3).Value = "CO2PAN"
$ws1.Cells.Item(29, 4).Value = 0.1437
$ws1.Cells.Item(29, 5).Value = "MT"

$ws1.Cells.Item(30, 1).Value = 1
$ws1.Cells.Item(30, 2).Value = "PWROILCRIXX01"
$ws1.Cells.Item(30, 3).Value = "CO2CRI"
$ws1.Cells.Item(30, 4).Value = 0.2021
$ws1.Cells.Item(30, 5).Value = "MT"

$ws1.Cells.Item(31, 1).Value = 2
$ws1.Cells.Item(31, 2).Value = "PWROILCRIXX01"
$ws1.Cells.Item(31, 3).Value = "CO2CRI"
$ws1.Cells.Item(31, 4).Value = 0.2021
$ws1.Cells.Item(31, 5).Value = "MT"

$ws1.Cells.Item(32, 1).Value = 1
$ws1.Cells.Item(32, 2).Value = "PWROILPANXX01"
$ws1.Cells.Item(32, 3).Value = "CO2PAN"
$ws1.Cells.Item(32, 4).Value = 0.2077
$ws1.Cells.Item(32, 5).Value = "MT"

$ws1.Cells.Item(33, 1).Value = 2
$ws1.Cells.Item(33, 2).Value = "PWROILPANXX01"
$ws1.Cells.Item(33, 3).Value = "CO2PAN"
$ws1.Cells.Item(33, 4).Value = 0.2077
$ws1.Cells.Item(33, 5).Value = "MT"

$ws1.Cells.Item(34, 1).Value = 1
$ws1.Cells.Item(34, 2).Value = "PWROTHCRIXX01"
$ws1.Cells.Item(34, 3).Value = "CO2CRI"
$ws1.Cells.Item(34, 4).Value = 0.0503
$ws1.Cells.Item(34, 5).Value = "MT"

$ws1.Cells.Item(35, 1).Value = 2
$ws1.Cells.Item(35, 2).Value = "PWROTHCRIXX01"
$ws1.Cells.Item(35, 3).Value = "CO2CRI"
$ws1.Cells.Item(35, 4).Value = 0.0503
$ws1.Cells.Item(35, 5).Value = "MT"

$ws1.Cells.Item(36, 1).Value = 1
$ws1.Cells.Item(36, 2).Value = "PWROTHPANXX01"
$ws1.Cells.Item(36, 3).Value = "CO2PAN"
$ws1.Cells.Item(36, 4).Value = 0.0503
$ws1.Cells.Item(36, 5).Value = "MT"

$ws1.Cells.Item(37, 1).Value = 2
$ws1.Cells.Item(37, 2).Value = "PWROTHPANXX01"
$ws1.Cells.Item(37, 3).Value = "CO2PAN"
$ws1.Cells.Item(37, 4).Value = 0.0503
$ws1.Cells.Item(37, 5).Value = "MT"

$ws1.Cells.Item(38, 1).Value = 1
$ws1.Cells.Item(38, 2).Value = "PWRPETCRIXX01"
$ws1.Cells.Item(38, 3).Value = "CO2CRI"
$ws1.Cells.Item(38, 4).Value = 0.217
$ws1.Cells.Item(38, 5).Value = "MT"

$ws1.Cells.Item(39, 1).Value = 2
$ws1.Cells.Item(39, 2).Value = "PWRPETCRIXX01"
$ws1.Cells.Item(39, 3).Value = "CO2CRI"
$ws1.Cells.Item(39, 4).Value = 0.217
$ws1.Cells.Item(39, 5).Value = "MT"

$ws1.Cells.Item(40, 1).Value = 1
$ws1.Cells.Item(40, 2).Value = "PWRPETPANXX01"
$ws1.Cells.Item(40, 3).Value = "CO2PAN"
$ws1.Cells.Item(40, 4).Value = 0.217
$ws1.Cells.Item(40, 5).Value = "MT"

$ws1.Cells.Item(41, 1).Value = 2
$ws1.Cells.Item(41, 2).Value = "PWRPETPANXX01"
$ws1.Cells.Item(41, 3).Value = "CO2PAN"
$ws1.Cells.Item(41, 4).Value = 0.217
$ws1.Cells.Item(41, 5).Value = "MT"

$ws1.Cells.Item(42, 1).Value = 1
$ws1.Cells.Item(42, 2).Value = "PWRWASCRIXX01"
$ws1.Cells.Item(42, 3).Value = "CO2CRI"
$ws1.Cells.Item(42, 4).Value = 0.0879
$ws1.Cells.Item(42, 5).Value = "MT"

$ws1.Cells.Item(43, 1).Value = 1
$ws1.Cells.Item(43, 2).Value = "PWRWASPANXX01"
$ws1.Cells.Item(43, 3).Value = "CO2PAN"
$ws1.Cells.Item(43, 4).Value = 0.0879
$ws1.Cells.Item(43, 5).Value = "MT"

# Externalities sheet: remove the now-obsolete CO2CRI externality row (row 2),
# shrinking the used range back down to just the header row.
$ws2 = $wb.Worksheets.Item(2)
$ws2.Rows.Item(2).Delete()
